# Ledger export change:
#  - Insert a new 3-column mini summary block (TITLE/DATE/TIME + one data row)
#    above the existing ledger table, pushing the ledger down by 3 rows
#    (original header row 1 -> row 4, data rows 2-7 -> rows 5-10). Row 3
#    stays empty, acting as a spacer between the two tables.
#  - The new header cells reuse the same bold header style as the ledger
#    header, and the new date/time values get dedicated number formats
#    (custom "yyyy-mm-dd" and the built-in "h:mm:ss").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing table down by 3 rows, leaving rows 1-3 free for the
# new summary block.
$ws.Rows("1:3").Insert()

# Reuse the existing bold "header" cell style (now on row 4) for the new
# mini-table header on row 1, instead of building a brand new font/style.
$ws.Range("A4:C4").Copy()
$ws.Range("A1:C1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A1").Value = "TITLE"
$ws.Range("B1").Value = "DATE"
$ws.Range("C1").Value = "TIME"

# New summary data row.
$ws.Range("A2").Value = "Cash To Mustafa"

$ws.Range("B2").Value = 45206
$ws.Range("B2").NumberFormat = "yyyy-mm-dd"

$ws.Range("C2").Value = 0.7096339854976852
$ws.Range("C2").NumberFormat = "h:mm:ss"
